$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.316.90"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -5.33%  "
$ws.Range("D3").Value = "'3.006.51"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.55%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'575.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").Value = "'126.22"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.66%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'3.001.42"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.66%  "
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("E10").Value = "  -7.99%  "
$ws.Range("E11").Value = "  -5.95%  "
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("D13").Value = "'0.0000222"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -7.90%  "
$ws.Range("D14").Value = "'32.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.50%  "
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "'3.496.83"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.61%  "
$ws.Range("D17").Value = "'3.001.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.59%  "
$ws.Range("D18").Value = "'60.198.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.48%  "
$ws.Range("D19").Value = "'6.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").Value = "'429.88"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.14%  "
$ws.Range("D21").Value = "'13.15"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").Value = "'0.670"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.62%  "
$ws.Range("D23").Value = "'7.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -8.11%  "
$ws.Range("D24").Value = "'12.94"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("D25").Value = "'79.45"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.78%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'2.55"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.05%  "
$ws.Range("D29").Value = "'1.97"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.62%  "
$ws.Range("D30").Value = "'7.30"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.56%  "
$ws.Range("E31").Value = "  -10.93%  "
$ws.Range("D32").Value = "'25.34"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.98%  "
$ws.Range("D33").Value = "'0.0945"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.98%  "
$ws.Range("E34").Value = "  -5.13%  "
$ws.Range("E35").Value = "  -9.11%  "
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("E37").Value = "  -15.24%  "
$ws.Range("D38").Value = "'0.0₃0677"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -8.90%  "
$ws.Range("D39").Value = "'8.48"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.89%  "
$ws.Range("D40").Value = "'0.0357"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -9.36%  "
$ws.Range("E41").Value = "  -5.09%  "
$ws.Range("D42").Value = "'377.15"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.54%  "
$ws.Range("D43").Value = "'2.677.71"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("D44").Value = "'2.47"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.68%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "'0.236"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.28%  "
$ws.Range("D47").Value = "'120.96"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.21%  "
$ws.Range("D48").Value = "'2.02"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.78%  "
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("D50").Value = "'23.63"
$ws.Range("D50").ClearFormats()
$ws.Range("E51").Value = "  -7.07%  "
